# "updated edit photos screen"
#
# 1) Re-sort the grants table on "Sheet1": the "photos" resource grants
#    move up (right after gallery_categories), "members" moves to the
#    bottom, and "posts" follows after it. The header in A3 is updated.
# 2) Add a new "photos" worksheet at the end of the workbook containing
#    the list of uploaded gallery photos and a generated SQL INSERT
#    statement per row.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Update the title in A3 -------------------------------------------------
$ws1.Range("A3").Value = "photographer: role id = 4"

# --- Re-sort the grants rows ------------------------------------------------
# Current layout (rows 12-20):
#   12            members / READ
#   13-16         photos  / CREATE,DELETE,READ,UPDATE
#   17-20 (s="6") posts   / CREATE,DELETE,READ,UPDATE
#
# New layout (rows 12-24):
#   12-15         photos  / CREATE,DELETE,READ,UPDATE
#   20            members / READ
#   21-24 (s="6") posts   / CREATE,DELETE,READ,UPDATE

# photos grants move up to rows 12-15
$ws1.Range("B12").Value = 74
$ws1.Range("C12").Value = "photos"
$ws1.Range("D12").Value = "CREATE"

$ws1.Range("B13").Value = 77
$ws1.Range("C13").Value = "photos"
$ws1.Range("D13").Value = "DELETE"

$ws1.Range("B14").Value = 76
$ws1.Range("C14").Value = "photos"
$ws1.Range("D14").Value = "READ"

$ws1.Range("B15").Value = 75
$ws1.Range("C15").Value = "photos"
$ws1.Range("D15").Value = "UPDATE"

# clear the now-stale rows 16-19 entirely (old posts rows used to live
# here with highlighted formatting; they move down to 21-24 below)
$ws1.Range("A16:D19").Clear()

# members grant moves down to row 20 (no special fill/border)
$ws1.Range("B20").Value = 25
$ws1.Range("C20").Value = "members"
$ws1.Range("D20").Value = "READ"

# posts grants move down to rows 21-24, keeping their highlighted style
$ws1.Range("B21").Value = 46
$ws1.Range("C21").Value = "posts"
$ws1.Range("D21").Value = "CREATE"

$ws1.Range("B22").Value = 49
$ws1.Range("C22").Value = "posts"
$ws1.Range("D22").Value = "DELETE"

$ws1.Range("B23").Value = 47
$ws1.Range("C23").Value = "posts"
$ws1.Range("D23").Value = "READ"

$ws1.Range("B24").Value = 48
$ws1.Range("C24").Value = "posts"
$ws1.Range("D24").Value = "UPDATE"

# Re-apply the highlighted (fill+border) formatting used by the
# "posts" block onto its new location (rows 21-24), taken from the
# still-highlighted "galleries" block in rows 4-7.
$ws1.Range("B4:D7").Copy() | Out-Null
$ws1.Range("B21:D24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# restore the selection/active cell as it ends up after the edits
$ws1.Range("A4").Select() | Out-Null

# --- Add the new "photos" worksheet at the end ------------------------------
$sheetCount = $wb.Worksheets.Count
$wsPhotos = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$wsPhotos.Name = "photos"

# header row
$wsPhotos.Range("A1").Value = "id"
$wsPhotos.Range("B1").Value = "gallery_id"
$wsPhotos.Range("C1").Value = "photographer_id"
$wsPhotos.Range("D1").Value = "photo_title"
$wsPhotos.Range("E1").Value = "description"
$wsPhotos.Range("F1").Value = "filename"
$wsPhotos.Range("G1").Value = "order"
$wsPhotos.Range("H1").Value = "created_at"
$wsPhotos.Range("I1").Value = "updated_at"
$wsPhotos.Range("J1").Value = "path"

# data rows: filenames / order, plus the generated INSERT statement
for ($i = 1; $i -le 26; $i++) {
    $row = $i + 1
    $filename = "GTU_Eggdrop_2012_{0:D2}.jpg" -f $i

    $wsPhotos.Range("F$row").Value = $filename
    $wsPhotos.Range("G$row").Value = $i

    $formula = '="INSERT INTO ``photos`` (``gallery_id``,``photographer_id``,``filename``,``order``,``created_at``,``updated_at``,``path``) VALUES (14,9,''"&F' + $row + '&"'',''"&G' + $row + '&"'',''2012-03-30'',''2012-03030'',''2012_egg_drop'');"'
    $wsPhotos.Range("K$row").Formula = $formula
}

# cosmetic: column widths that "best fit" the content
$wsPhotos.Columns.Item(2).ColumnWidth = 23.5
$wsPhotos.Columns.Item(3).ColumnWidth = 15
$wsPhotos.Columns.Item(4).ColumnWidth = 10.33
$wsPhotos.Columns.Item(5).ColumnWidth = 10.33
$wsPhotos.Columns.Item(6).ColumnWidth = 23.5

$wsPhotos.Range("K2:K27").Select() | Out-Null

# keep "Sheet1" as the selected/active tab, as before the edit
$ws1.Activate() | Out-Null
